# Updates cryptocurrency listing data (Price, Volume label, Hora) to the
# "31-12-2022 12:xx" refresh snapshot, per the Dec 31 2022 GitHub Actions run.
#
# All target columns (B, C, D, E, G) hold text in the source workbook (every
# cell is serialized as an inline string, even the numeric-looking Price/Hora
# columns), so values are written with a leading apostrophe to force Excel to
# keep them as text instead of silently re-typing them as numbers. The
# apostrophe triggers Excel's "number stored as text" quote-prefix cell
# style; resetting the range back to the "Normal" style afterwards clears
# that cosmetic style flag while leaving the text value/type intact, so the
# written cells land exactly like the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}


Set-TextValue "D2" "245.66"
Set-TextValue "G2" "12"
Set-TextValue "D3" "26.35"
Set-TextValue "G3" "12"
Set-TextValue "D4" "5.102"
Set-TextValue "G4" "12"
Set-TextValue "D5" "0.05596"
Set-TextValue "G5" "12"
Set-TextValue "D6" "6.481"
Set-TextValue "G6" "12"
Set-TextValue "D7" "3.035"
Set-TextValue "G7" "12"
Set-TextValue "D8" "0.8114"
Set-TextValue "G8" "12"
Set-TextValue "D9" "0.8414"
Set-TextValue "G9" "12"
Set-TextValue "D10" "0.1346"
Set-TextValue "G10" "12"
Set-TextValue "D11" "0.02840"
Set-TextValue "G11" "12"
Set-TextValue "D12" "0.09401"
Set-TextValue "G12" "12"
Set-TextValue "D13" "0.001516"
Set-TextValue "G13" "12"
Set-TextValue "D14" "0.0006011"
Set-TextValue "E14" "13OneONEWorstin24h"
Set-TextValue "G14" "12"
Set-TextValue "D15" "0.006115"
Set-TextValue "G15" "12"
Set-TextValue "D16" "3.564"
Set-TextValue "G16" "12"
Set-TextValue "G17" "12"
Set-TextValue "G18" "12"
Set-TextValue "D19" "0.07014"
Set-TextValue "G19" "12"
Set-TextValue "G20" "12"
Set-TextValue "G21" "12"
Set-TextValue "D22" "3.741"
Set-TextValue "G22" "12"
Set-TextValue "D23" "0.04687"
Set-TextValue "G23" "12"
Set-TextValue "G24" "12"
Set-TextValue "D25" "0.001248"
Set-TextValue "G25" "12"
Set-TextValue "D26" "0.004598"
Set-TextValue "G26" "12"
Set-TextValue "D27" "0.00009600"
Set-TextValue "G27" "12"
Set-TextValue "G28" "12"
Set-TextValue "G29" "12"
Set-TextValue "G30" "12"
Set-TextValue "G31" "12"
Set-TextValue "G32" "12"
Set-TextValue "G33" "12"
Set-TextValue "G34" "12"
Set-TextValue "G35" "12"
Set-TextValue "G36" "12"
Set-TextValue "G37" "12"
Set-TextValue "G38" "12"
Set-TextValue "G39" "12"
Set-TextValue "D40" "0.03654"
Set-TextValue "G40" "12"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1367"
Set-TextValue "E41" "40BKEXTokenBKKBestin24h"
Set-TextValue "G41" "12"
Set-TextValue "B42" "CEJI"
Set-TextValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002660"
Set-TextValue "E42" "41CEJICEJI"
Set-TextValue "G42" "12"
Set-TextValue "B43" "KickToken"
Set-TextValue "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003403"
Set-TextValue "E43" "42KickTokenKICK"
Set-TextValue "G43" "12"
Set-TextValue "D44" "0.008549"
Set-TextValue "G44" "12"
Set-TextValue "D45" "0.00005296"
Set-TextValue "G45" "12"
Set-TextValue "G46" "12"
Set-TextValue "E47" "46CoinbaseStockTokenCOIN"
Set-TextValue "G47" "12"
Set-TextValue "G48" "12"
Set-TextValue "G49" "12"
Set-TextValue "G50" "12"
Set-TextValue "G51" "12"
